$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.533.68"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "3.628.96"
$ws.Range("E3").Value = "  +2.90%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "195.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  +5.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.648"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.39"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000304"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("D14").Value = "4.205.59"
$ws.Range("E14").Value = "  +2.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "603.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "12.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.05%  "
$ws.Range("D17").Value = "70.631.96"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.629.65"
$ws.Range("E18").Value = "  +3.29%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("E20").Value = "  +1.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  -4.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.30"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.54%  "
$ws.Range("E33").Value = "  +2.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("D35").Value = "0.0₃0887"
$ws.Range("E35").Value = "  +2.41%  "
$ws.Range("D36").Value = "3.918.09"
$ws.Range("E36").Value = "  +5.28%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "521.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.81%  "
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.389"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.136"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0464"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.71%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000250"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.35%  "
$ws.Range("E51").Value = "  +1.93%  "
